# Updated symbol list on Sat Dec 17 07:27:38 UTC 2022 with GitHub Actions
#
# All values in this sheet are stored as plain text (inline strings), even
# when they look like numbers (e.g. "230.76"). Writing a numeric-looking
# string straight into .Value/.Value2 makes Excel coerce the cell to a
# real number (and can introduce floating point noise, e.g.
# 230.74000000000001). To keep the cells as exact text - matching the
# original sheet - force the cell to Text format before the write, then
# restore whatever style the cell had before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $range = $ws.Range($Address)
    $oldStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value2 = $Value
    $range.Style = $oldStyle
}

# Price (column D) corrections
Set-TextValue "D2"  "230.74"
Set-TextValue "D3"  "22.60"
Set-TextValue "D4"  "5.293"
Set-TextValue "D6"  "3.382"
Set-TextValue "D7"  "6.470"
Set-TextValue "D8"  "1.062"
Set-TextValue "D9"  "0.7821"
Set-TextValue "D10" "0.1397"
Set-TextValue "D11" "0.07410"
Set-TextValue "D12" "0.03153"
Set-TextValue "D13" "0.02969"
Set-TextValue "D15" "0.001660"
Set-TextValue "D16" "3.257"
Set-TextValue "D17" "0.04752"
Set-TextValue "D18" "0.0005793"
Set-TextValue "D19" "0.006230"
Set-TextValue "D20" "0.005225"
Set-TextValue "D23" "3.970"

# Row 27 (UpBots) - price + volume label tweak
Set-TextValue "D27" "0.0004993"
Set-TextValue "E27" "26UpBotsUBXTBestin24h"

Set-TextValue "D40" "0.04034"

# Row 41 (KickToken) - price + volume label tweak
Set-TextValue "D41" "0.006993"
Set-TextValue "E41" "40KickTokenKICK"

# Rows 42/43 swap places (BKEXToken <-> CEJI) with refreshed data
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1042"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002940"
Set-TextValue "E43" "42CEJICEJI"

Set-TextValue "D44" "0.009240"
Set-TextValue "D45" "0.00005443"
Set-TextValue "D47" "0.7857"
Set-TextValue "D48" "0.04096"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.01011"
